# The deck's slide-master theme ("Integral" / "Red Violet") is switched
# over to the stock PowerPoint "Office Theme" color palette - i.e. the
# same effect as opening Design > Variants > Colors and picking the
# built-in "Office" swatch. Font scheme and effect scheme were already
# the stock "Office" ones, so only the 12 theme colors need updating.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster

# Classic 8-slot ColorScheme: Colors(1..8) = dk1, lt1, dk2, lt2,
# accent1, accent2, accent3, accent4. Values are VBA RGB()-packed
# integers (R + G*256 + B*65536).
$cs = $master.ColorScheme
$cs.Colors(1).RGB = 0         # dk1      000000
$cs.Colors(2).RGB = 16777215  # lt1      FFFFFF
$cs.Colors(3).RGB = 6968388   # dk2      44546A
$cs.Colors(4).RGB = 15132391  # lt2      E7E6E6
$cs.Colors(5).RGB = 13998939  # accent1  5B9BD5
$cs.Colors(6).RGB = 3243501   # accent2  ED7D31
$cs.Colors(7).RGB = 10855845  # accent3  A5A5A5
$cs.Colors(8).RGB = 49407     # accent4  FFC000

# Full 12-slot DrawingML theme color scheme covers the remaining
# accent5, accent6, hlink and folHlink entries (items 9-12) that the
# legacy ColorScheme object doesn't expose.
$tcs = $p.Slides.Item(1).ThemeColorScheme
$tcs.Item(9).RGB  = 12874308  # accent5  4472C4
$tcs.Item(10).RGB = 4697456   # accent6  70AD47
$tcs.Item(11).RGB = 12673797  # hlink    0563C1
$tcs.Item(12).RGB = 7491477   # folHlink 954F72
